$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new rows for t-value and p-value of the correlation coefficient
$ws.Range("B22").Value = "t値"
$ws.Range("C22").Formula = "=C20*SQRT((C21-2)/(1-C20^2))"

$ws.Range("B23").Value = "p値"
$ws.Range("C23").Formula = "=T.DIST.2T(ABS(C22), C21-2)"

# Update the selection to match the new content
$ws.Range("D23").Select()
